# Update rows 7, 8 and 9 on the active sheet of the "mapa_interactivo_Optical_Power"
# workbook with the latest claims data (new case added, others shifted up).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (Caso) and B (F. De Reclamo) hold numeric-looking / date-looking
# text in this sheet (e.g. "6100", "8/8/2025"), so force them to be stored
# as text instead of being auto-converted to a number/date.
$ws.Range("A7:B9").NumberFormat = "@"

# ---- Row 7 ----
$ws.Range("A7").Value = "6100"
$ws.Range("B7").Value = "8/8/2025"
$ws.Range("C7").Value = "DE LOS CONSTITUYENTES AV. 5552"
$ws.Range("D7").Value = 12
$ws.Range("H7").Value = "Cable en panza y cables cortados"
$ws.Range("J7").Value = '{"direccionesNormalizadas": [{"altura": 5552, "cod_calle": 4043, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.501174", "y": "-34.575005"}, "direccion": "DE LOS CONSTITUYENTES AV. 5552, CABA", "nombre_calle": "DE LOS CONSTITUYENTES AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K7").Value = -58.501174
$ws.Range("L7").Value = -34.575005
$ws.Range("M7").Value = "Paternal"
$ws.Range("N7").Value = "Capital Norte"

# ---- Row 8 ----
$ws.Range("A8").Value = "6265"
$ws.Range("B8").Value = "8/7/2025"
$ws.Range("C8").Value = "BROWN, ALTE. AV. 881"
$ws.Range("D8").Value = 4
$ws.Range("H8").Value = "tendido a baja altura"
$ws.Range("J8").Value = '{"direccionesNormalizadas": [{"altura": 881, "cod_calle": 2115, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.360551", "y": "-34.632684"}, "direccion": "BROWN, ALTE. AV. 881, CABA", "nombre_calle": "BROWN, ALTE. AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K8").Value = -58.360551
$ws.Range("L8").Value = -34.632684
$ws.Range("M8").Value = "San Telmo"
$ws.Range("N8").Value = "Capital Sur"

# ---- Row 9 ----
$ws.Range("A9").Value = "6570"
$ws.Range("B9").Value = "8/8/2025"
$ws.Range("C9").Value = "FALCON, RAMON L.,CNEL. 7096"
$ws.Range("D9").Value = 9
$ws.Range("H9").Value = "Tendido a muy baja altura"
$ws.Range("J9").Value = '{"direccionesNormalizadas": [{"altura": 7096, "cod_calle": 6006, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.526204", "y": "-34.640276"}, "direccion": "FALCON, RAMON L.,CNEL. 7096, CABA", "nombre_calle": "FALCON, RAMON L.,CNEL.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K9").Value = -58.526204
$ws.Range("L9").Value = -34.640276
$ws.Range("M9").Value = "Devoto"
$ws.Range("N9").Value = "Capital Norte"
